# Insert a new weekly price record for "Albahaca" (Terminal La Palmera de
# La Serena) as row 125, shifting all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 125 (existing rows 125-172 shift to 126-173).
$ws.Rows.Item(125).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A125").Value2 = 8
$ws.Range("B125").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C125").Value2 = "Coquimbo"
$ws.Range("D125").Value2 = 45009
$ws.Range("E125").Value2 = 4
$ws.Range("F125").Value2 = 100112052
$ws.Range("G125").Value2 = "Albahaca"
$ws.Range("H125").Value2 = "Sin especificar"
$ws.Range("I125").Value2 = "Primera"
$ws.Range("J125").Value2 = 1060
$ws.Range("K125").Value2 = 2500
$ws.Range("L125").Value2 = 3000
$ws.Range("M125").Value2 = 2750
$ws.Range("N125").Value2 = "`$/docena de matas"
$ws.Range("O125").Value2 = "Provincia del Elquí"
$ws.Range("P125").Value2 = 458
$ws.Range("Q125").Value2 = 6
$ws.Range("R125").Value2 = "Hortaliza"
